# Edit script for cambodia_cpl_2023-2024.xlsx
# 1) Swap the two duplicate-date matches that were recorded in the wrong row order (rows 25 & 26)
# 2) Append 5 new match rows (42-46) scraped in the 31-10-2023 run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap F:V content between row 25 and row 26 ---
$row25 = $ws.Range("F25:V25").Value()
$row26 = $ws.Range("F26:V26").Value()
$ws.Range("F25:V25").Value = $row26
$ws.Range("F26:V26").Value = $row25

# --- Step 2: append new rows 42-46, copying the formatting of row 41 first ---
$ws.Range("A41:V41").Copy()
$ws.Range("A42:V42").PasteSpecial(-4122)
$ws.Range("A43:V43").PasteSpecial(-4122)
$ws.Range("A44:V44").PasteSpecial(-4122)
$ws.Range("A45:V45").PasteSpecial(-4122)
$ws.Range("A46:V46").PasteSpecial(-4122)

# Row 42
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 5).Value = 45227.54166666666
$ws.Cells.Item(42, 2).Value = "cambodia"
$ws.Cells.Item(42, 3).Value = "cpl"
$ws.Cells.Item(42, 4).Value = "2023-2024"
$ws.Cells.Item(42, 6).Value = "Boeung Ket"
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = "Svay Rieng"
$ws.Cells.Item(42, 9).Value = 3
$ws.Cells.Item(42, 10).Value = 2.93
$ws.Cells.Item(42, 11).Value = "27/10/2023 01:12"
$ws.Cells.Item(42, 12).Value = 3.18
$ws.Cells.Item(42, 13).Value = "28/10/2023 11:12"
$ws.Cells.Item(42, 14).Value = 3.44
$ws.Cells.Item(42, 15).Value = "27/10/2023 01:12"
$ws.Cells.Item(42, 16).Value = 3.86
$ws.Cells.Item(42, 17).Value = "28/10/2023 11:12"
$ws.Cells.Item(42, 18).Value = 1.97
$ws.Cells.Item(42, 19).Value = "27/10/2023 01:12"
$ws.Cells.Item(42, 20).Value = 1.88
$ws.Cells.Item(42, 21).Value = "28/10/2023 11:12"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/boeung-ket-svay-rieng/vNTVz71d/"

# Row 43
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 5).Value = 45228.40625
$ws.Cells.Item(43, 2).Value = "cambodia"
$ws.Cells.Item(43, 3).Value = "cpl"
$ws.Cells.Item(43, 4).Value = "2023-2024"
$ws.Cells.Item(43, 6).Value = "Angkor Tiger"
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = "NagaWorld"
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 3.7
$ws.Cells.Item(43, 11).Value = "27/10/2023 23:12"
$ws.Cells.Item(43, 12).Value = 4.56
$ws.Cells.Item(43, 13).Value = "29/10/2023 09:30"
$ws.Cells.Item(43, 14).Value = 3.71
$ws.Cells.Item(43, 15).Value = "27/10/2023 23:12"
$ws.Cells.Item(43, 16).Value = 4.03
$ws.Cells.Item(43, 17).Value = "29/10/2023 09:30"
$ws.Cells.Item(43, 18).Value = 1.64
$ws.Cells.Item(43, 19).Value = "27/10/2023 23:12"
$ws.Cells.Item(43, 20).Value = 1.56
$ws.Cells.Item(43, 21).Value = "29/10/2023 09:30"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/angkor-tiger-nagaworld/I5YAbAg3/"

# Row 44
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 5).Value = 45228.40625
$ws.Cells.Item(44, 2).Value = "cambodia"
$ws.Cells.Item(44, 3).Value = "cpl"
$ws.Cells.Item(44, 4).Value = "2023-2024"
$ws.Cells.Item(44, 6).Value = "Kirivong Sok Sen Chey"
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = "Phnom Penh Crown"
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 6.09
$ws.Cells.Item(44, 11).Value = "26/10/2023 23:12"
$ws.Cells.Item(44, 12).Value = 8.8
$ws.Cells.Item(44, 13).Value = "29/10/2023 09:30"
$ws.Cells.Item(44, 14).Value = 4.94
$ws.Cells.Item(44, 15).Value = "26/10/2023 23:12"
$ws.Cells.Item(44, 16).Value = 5.68
$ws.Cells.Item(44, 17).Value = "29/10/2023 09:30"
$ws.Cells.Item(44, 18).Value = 1.29
$ws.Cells.Item(44, 19).Value = "26/10/2023 23:12"
$ws.Cells.Item(44, 20).Value = 1.22
$ws.Cells.Item(44, 21).Value = "29/10/2023 09:30"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/kirivong-sok-sen-chey-phnom-penh-crown/0WURymnj/"

# Row 45
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 5).Value = 45228.40625
$ws.Cells.Item(45, 2).Value = "cambodia"
$ws.Cells.Item(45, 3).Value = "cpl"
$ws.Cells.Item(45, 4).Value = "2023-2024"
$ws.Cells.Item(45, 6).Value = "Prey Veng"
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = "Visakha"
$ws.Cells.Item(45, 9).Value = 2
$ws.Cells.Item(45, 10).Value = 4.07
$ws.Cells.Item(45, 11).Value = "27/10/2023 23:12"
$ws.Cells.Item(45, 12).Value = 6.02
$ws.Cells.Item(45, 13).Value = "29/10/2023 09:30"
$ws.Cells.Item(45, 14).Value = 4.12
$ws.Cells.Item(45, 15).Value = "27/10/2023 23:12"
$ws.Cells.Item(45, 16).Value = 4.95
$ws.Cells.Item(45, 17).Value = "29/10/2023 09:30"
$ws.Cells.Item(45, 18).Value = 1.51
$ws.Cells.Item(45, 19).Value = "27/10/2023 23:12"
$ws.Cells.Item(45, 20).Value = 1.36
$ws.Cells.Item(45, 21).Value = "29/10/2023 09:30"
$ws.Cells.Item(45, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/prey-veng-visakha/CSXEcU89/"

# Row 46
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 5).Value = 45228.5
$ws.Cells.Item(46, 2).Value = "cambodia"
$ws.Cells.Item(46, 3).Value = "cpl"
$ws.Cells.Item(46, 4).Value = "2023-2024"
$ws.Cells.Item(46, 6).Value = "Dangkor"
$ws.Cells.Item(46, 7).Value = 2
$ws.Cells.Item(46, 8).Value = "Tiffy Army"
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 4.74
$ws.Cells.Item(46, 11).Value = "28/10/2023 01:12"
$ws.Cells.Item(46, 12).Value = 3.66
$ws.Cells.Item(46, 13).Value = "29/10/2023 11:46"
$ws.Cells.Item(46, 14).Value = 4.33
$ws.Cells.Item(46, 15).Value = "28/10/2023 01:12"
$ws.Cells.Item(46, 16).Value = 4.18
$ws.Cells.Item(46, 17).Value = "29/10/2023 11:46"
$ws.Cells.Item(46, 18).Value = 1.43
$ws.Cells.Item(46, 19).Value = "28/10/2023 01:12"
$ws.Cells.Item(46, 20).Value = 1.68
$ws.Cells.Item(46, 21).Value = "29/10/2023 11:46"
$ws.Cells.Item(46, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/dangkor-senchey-tiffy-army/vFZ6ajvc/"
